# Updated cryptos list on Mon Aug  7 17:32:35 UTC 2023 with GitHub Actions
#
# Refreshes the "Price" (column D) and "Volume(1h)" (column E) figures for
# every coin row (2-51) to the latest scraped values. Row 49/50 additionally
# swap rank position between "TheSandbox" and "EnergySwap" (name + link are
# updated along with their price/volume).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D values are decimal-looking strings (e.g. "1.000", "0.000008920",
# "28.939.69") that must be preserved exactly as text. If we simply assign a
# plain numeric-looking string to Range.Value, Excel auto-converts it to a
# floating point number and mangles the text (drops trailing zeros, adds
# binary rounding noise, etc.), so the whole column is forced to the "@"
# (Text) number format first and restored to the default formatting again
# once every value has been written.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "28.939.69"
$ws.Range("E2").Value = "  -0.38%  "

$ws.Range("D3").Value = "1.812.87"
$ws.Range("E3").Value = "  -0.99%  "

$ws.Range("D4").Value = "0.9994"
$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").Value = "241.16"
$ws.Range("E5").Value = "  -0.96%  "

$ws.Range("D6").Value = "0.6095"
$ws.Range("E6").Value = "  -2.97%  "

$ws.Range("D7").Value = "1.002"
$ws.Range("E7").Value = "  +0.11%  "

$ws.Range("D8").Value = "0.07291"
$ws.Range("E8").Value = "  -2.56%  "

$ws.Range("D9").Value = "0.2868"
$ws.Range("E9").Value = "  -1.88%  "

$ws.Range("D10").Value = "22.78"
$ws.Range("E10").Value = "  -1.58%  "

$ws.Range("D11").Value = "0.07634"
$ws.Range("E11").Value = "  -1.07%  "

$ws.Range("D12").Value = "1.819.95"
$ws.Range("E12").Value = "  -0.51%  "

$ws.Range("D13").Value = "4.924"
$ws.Range("E13").Value = "  -1.34%  "

$ws.Range("D14").Value = "0.6588"
$ws.Range("E14").Value = "  -1.32%  "

$ws.Range("D15").Value = "80.99"
$ws.Range("E15").Value = "  -1.88%  "

$ws.Range("D16").Value = "0.000008920"
$ws.Range("E16").Value = "  -4.20%  "

$ws.Range("E17").Value = "  -2.72%  "

$ws.Range("D18").Value = "28.938.27"
$ws.Range("E18").Value = "  -0.50%  "

$ws.Range("D19").Value = "2.063.52"
$ws.Range("E19").Value = "  -0.97%  "

$ws.Range("D20").Value = "236.22"
$ws.Range("E20").Value = "  +5.95%  "

$ws.Range("D21").Value = "12.39"
$ws.Range("E21").Value = "  -1.44%  "

$ws.Range("E22").Value = "  -0.06%  "

$ws.Range("D23").Value = "7.087"
$ws.Range("E23").Value = "  -0.43%  "

$ws.Range("D24").Value = "1.000"
$ws.Range("E24").Value = "  -0.01%  "

$ws.Range("D25").Value = "158.40"
$ws.Range("E25").Value = "  -0.66%  "

$ws.Range("D26").Value = "0.1400"
$ws.Range("E26").Value = "  +0.16%  "

$ws.Range("D27").Value = "8.390"
$ws.Range("E27").Value = "  -1.28%  "

$ws.Range("E28").Value = "  -2.04%  "

$ws.Range("D29").Value = "1.479"
$ws.Range("E29").Value = "  -1.23%  "

$ws.Range("D30").Value = "0.05587"
$ws.Range("E30").Value = "  -2.70%  "

$ws.Range("D31").Value = "4.059"
$ws.Range("E31").Value = "  +0.00%  "

$ws.Range("D32").Value = "4.069"
$ws.Range("E32").Value = "  -1.94%  "

$ws.Range("D33").Value = "1.209"
$ws.Range("E33").Value = "  +0.38%  "

$ws.Range("D34").Value = "1.815"
$ws.Range("E34").Value = "  -1.47%  "

$ws.Range("D35").Value = "0.7307"
$ws.Range("E35").Value = "  -2.44%  "

$ws.Range("E36").Value = "  -0.83%  "

$ws.Range("D37").Value = "2.630"
$ws.Range("E37").Value = "  -1.49%  "

$ws.Range("D38").Value = "2.803"
$ws.Range("E38").Value = "  +1.51%  "

$ws.Range("E39").Value = "  -1.90%  "

$ws.Range("D40").Value = "1.192.17"
$ws.Range("E40").Value = "  -2.17%  "

$ws.Range("D41").Value = "6.329"
$ws.Range("E41").Value = "  -3.29%  "

$ws.Range("D42").Value = "0.8913"
$ws.Range("E42").Value = "  -0.14%  "

$ws.Range("E43").Value = "  +0.01%  "

$ws.Range("D44").Value = "100.36"
$ws.Range("E44").Value = "  -1.79%  "

$ws.Range("D45").Value = "1.969.00"
$ws.Range("E45").Value = "  -0.64%  "

$ws.Range("D46").Value = "64.11"
$ws.Range("E46").Value = "  -2.42%  "

$ws.Range("D47").Value = "0.5079"
$ws.Range("E47").Value = "  -0.10%  "

$ws.Range("D48").Value = "0.00000000120"
$ws.Range("E48").Value = "  -5.76%  "

# Rows 49/50 swap rank: EnergySwap moves up to row 49, TheSandbox to row 50.
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "9.013"
$ws.Range("E49").Value = "  -0.39%  "

$ws.Range("B50").Value = "TheSandbox"
$ws.Range("C50").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D50").Value = "0.3971"
$ws.Range("E50").Value = "  -2.39%  "

$ws.Range("D51").Value = "0.05773"
$ws.Range("E51").Value = "  -0.84%  "

# Restore the default (General) number formatting on column D now that every
# value has been written, so cells don't carry a residual custom format.
$ws.Range("D2:D51").ClearFormats()
